$wb = $excel.ActiveWorkbook

# --- Sheet "mock": numeric value updates ---
$mock = $wb.Worksheets.Item("mock")

# Row 2 (Homo sapiens)
$mock.Range("C2").Value = 87130
$mock.Range("D2").Value = 60576
$mock.Range("E2").Value = 87130
$mock.Range("F2").Value = 34171
$mock.Range("G2").Value = 28324
$mock.Range("H2").Value = 18270
$mock.Range("I2").Value = 20396

# Row 6 (Zea mays)
$mock.Range("C6").Value = 50
$mock.Range("D6").Value = 50
$mock.Range("E6").Value = 50
$mock.Range("F6").Value = 4
$mock.Range("G6").Value = 4
$mock.Range("H6").Value = 4
$mock.Range("I6").Value = 4

# Row 12 (Lactobacillus)
$mock.Range("C12").Value = 4

# Row 13 (Lactobacillus helveticus)
$mock.Range("E13").Value = 4

# Row 19 (Methylobacterium radiotolerans)
$mock.Range("F19").Value = 5
$mock.Range("G19").Value = 5
$mock.Range("H19").Value = 5
$mock.Range("I19").Value = 5

# Reset the "mock" sheet view: drop the scrolled-down top-left cell and the
# A27:I30 selection, returning to the default top-of-sheet A1 selection
# while keeping this sheet as the active tab.
$mock.Activate()
$mock.Range("A1").Select()

# --- Sheet "notes": update the collapsed-taxa rule text and selection ---
$notes = $wb.Worksheets.Item("notes")
$notes.Range("A8").Value = "collapsed < mintaxa (5)"

# Move the notes sheet selection from A10 to A9, then restore "mock" as the
# active tab (selecting on "notes" would otherwise leave it active).
$notes.Range("A9").Select()
$mock.Activate()
